# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1=Wins, AE1=Losses, AF1=Ties ---
# Copy formatting from the neighboring header cell (AC1) so the new
# header cells get the same bold/border/centered style, then overwrite
# the copied text with the correct header labels.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-39): team record repeated for every player row ---
$lastRow = 39
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 91   # column AD - Wins
    $ws.Cells.Item($r, 31).Value = 71   # column AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # column AF - Ties
}
